# Generate Report for Handoff
# Adds two newly-handed-off files (2332db29... and 99c3d242...) to the
# Overview / zh-cn / de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$commit = "beafd4dfd5f87ad01dec48bcc2c2e88ba96fd99c"
$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/"

$hyperColor = 15570276   # OLE (BGR) form of RGB(0x64,0x95,0xED) == FF6495ED
$dateFmt = "yyyy-mm-dd HH:mm:ss"
$autoWidth = 16.38        # closest ColumnWidth achievable to the autosized 17.216 char-width

function Style-AsHyperlink($cell) {
    $cell.Font.Color = $hyperColor
    $cell.Font.Underline = $true
}

function Style-AsDate($cell) {
    $cell.NumberFormat = $dateFmt
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$newOverviewRows = @(
    @{ Row = 4; Name = "2332db29-2ae3-44e5-b9b0-4aed0d29a3f6.md" },
    @{ Row = 5; Name = "99c3d242-352d-4d91-bdfa-ff54b2f604b1.md" }
)

foreach ($item in $newOverviewRows) {
    $r = $item.Row
    $name = $item.Name
    $display = "e2e\" + $name

    $wsOverview.Cells.Item($r, 1).Value = $name
    $wsOverview.Cells.Item($r, 2).Value = $display
    $wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($r, 2), ($repoBase + $name), "", "", $display)
    Style-AsHyperlink $wsOverview.Cells.Item($r, 2)
    $wsOverview.Cells.Item($r, 3).Value = ".md"
    $wsOverview.Cells.Item($r, 4).Value = ""
    $wsOverview.Cells.Item($r, 5).Value = "Ready for handoff"
    $wsOverview.Cells.Item($r, 6).Value = "Ready for handoff"
    $wsOverview.Cells.Item($r, 7).Value = "2016-10-13 12:56:04"
    Style-AsDate $wsOverview.Cells.Item($r, 7)
}

$wsOverview.Columns.Item(5).ColumnWidth = $autoWidth
$wsOverview.Columns.Item(6).ColumnWidth = $autoWidth

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------
# Locale sheets (zh-cn / de-de) share the same column layout:
# A Source File Name (hyperlink) | B File Extension | C Status | D Source Path
# E Priority | F Content Duplicate | G Latest Handoff File | H Latest Handoff Datetime
# I Latest Target File | J Latest Handback File | K Latest Handback DateTime
# L Reference Tokens | M To be localized | N Dependency From | O Has metadata | P Error Detail
# ---------------------------------------------------------------------
$localeSheets = @(
    @{
        Name = "zh-cn";
        XlfSuffix = "zh-cn";
        HandoffTime = "2016-10-13 12:55:54"
    },
    @{
        Name = "de-de";
        XlfSuffix = "de-de";
        HandoffTime = "2016-10-13 12:56:04"
    }
)

$newLocaleFiles = @(
    @{ Name = "2332db29-2ae3-44e5-b9b0-4aed0d29a3f6.md"; Stem = "2332db29-2ae3-44e5-b9b0-4aed0d29a3f6"; Hash = "6c28bba9764a2a94c2025989009275e5d48fe179" },
    @{ Name = "99c3d242-352d-4d91-bdfa-ff54b2f604b1.md"; Stem = "99c3d242-352d-4d91-bdfa-ff54b2f604b1"; Hash = "0d3d472afef0fe334974fec3b7c2cb4b43179c29" }
)

foreach ($locale in $localeSheets) {
    $ws = $wb.Worksheets.Item($locale.Name)
    $row = 4
    foreach ($file in $newLocaleFiles) {
        $xlfName = $file.Stem + "." + $file.Hash + "." + $locale.XlfSuffix + ".xlf"

        $ws.Cells.Item($row, 1).Value = $file.Name
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 1), ($repoBase + $file.Name), "", "", $file.Name)
        Style-AsHyperlink $ws.Cells.Item($row, 1)

        $ws.Cells.Item($row, 2).Value = ".md"
        $ws.Cells.Item($row, 3).Value = "Ready for handoff"
        $ws.Cells.Item($row, 4).Value = "e2e"
        $ws.Cells.Item($row, 5).Value = "ht"
        # leading apostrophe forces text storage instead of Excel auto-coercing
        # the literal word True/False into a genuine boolean cell
        $ws.Cells.Item($row, 6).Value = "'False"
        $ws.Cells.Item($row, 7).Value = $xlfName
        $ws.Cells.Item($row, 8).Value = $locale.HandoffTime
        Style-AsDate $ws.Cells.Item($row, 8)
        $ws.Cells.Item($row, 11).Value = "0001-01-01 00:00:00"
        Style-AsDate $ws.Cells.Item($row, 11)
        $ws.Cells.Item($row, 13).Value = "'True"
        $ws.Cells.Item($row, 15).Value = "'False"

        $row = $row + 1
    }

    $ws.Columns.Item(3).ColumnWidth = $autoWidth

    $lo = $ws.ListObjects.Item(1)
    $lo.Resize($ws.Range("A1:P5"))
}

Write-Output "Added handoff rows for 2332db29-2ae3-44e5-b9b0-4aed0d29a3f6.md and 99c3d242-352d-4d91-bdfa-ff54b2f604b1.md"
